$d = $word.ActiveDocument

function Wrap-Xml([string]$bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

function Replace-ParaRuns([int]$paraIndex, [string]$bodyXml) {
    $p = $d.Paragraphs.Item($paraIndex).Range
    $full = $d.Range($p.Start, $p.End - 1)
    $full.InsertXML((Wrap-Xml $bodyXml))
}

$rPr = '<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'

# --- Edit 1: paragraph 20 -- "Expected behavior: If the items ..." split into 3 runs ---
$body1 = '<w:p>' +
    '<w:r>' + $rPr + '<w:t xml:space="preserve">Expected behavior: </w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t>Yes,ii</w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t>f the items are sufficient in inventory all the beverages will be prepared. A &#8220;beverage_name prepared&#8221; message should appear</w:t></w:r>' +
    '</w:p>'
Replace-ParaRuns 20 $body1

# --- Edit 2: paragraph 29 -- "Expected behavior: Yes. User can place a new request ..." split into 4 runs ---
$body2 = '<w:p>' +
    '<w:r>' + $rPr + '<w:t>Expected behavior: Yes. User can place a new request i</w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t xml:space="preserve">f </w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t>outlets are available</w:t></w:r>' +
    '</w:p>'
Replace-ParaRuns 29 $body2

# --- Edit 3: paragraph 31 -- add <w:lastRenderedPageBreak/> before the text ---
$body3 = '<w:p>' +
    '<w:r>' + $rPr + '<w:lastRenderedPageBreak/><w:t>Does the user see an error message for incorrect selection of beverages?</w:t></w:r>' +
    '</w:p>'
Replace-ParaRuns 31 $body3

# --- Edit 4: paragraph 41 -- "Is the user restricted ..." split into 3 runs ---
$body4 = '<w:p>' +
    '<w:r>' + $rPr + '<w:t>Is the user restricted to select only beverages up to</w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> N</w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> number of outlets per requests?</w:t></w:r>' +
    '</w:p>'
Replace-ParaRuns 41 $body4

# --- Edit 5: paragraph 45 -- append a new run with extra sentence ---
$body5 = '<w:p>' +
    '<w:r>' + $rPr + '<w:t>Expected behavior: No. The state of the machine is set as Refill. Only after refill the user can select beverages.</w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> The user will receive a message a to wait till refill is done.</w:t></w:r>' +
    '</w:p>'
Replace-ParaRuns 45 $body5

# --- Edit 6: move <w:lastRenderedPageBreak/> from paragraph 59 to paragraph 60 ---
$body6a = '<w:p>' +
    '<w:r>' + $rPr + '<w:t>After valid selection, does the machine show Beverages served?</w:t></w:r>' +
    '</w:p>'
Replace-ParaRuns 59 $body6a

$body6b = '<w:p>' +
    '<w:r>' + $rPr + '<w:lastRenderedPageBreak/><w:t>Expected behavior: Yes. A message is show for each beverage that is dispensed.</w:t></w:r>' +
    '</w:p>'
Replace-ParaRuns 60 $body6b

Write-Host "Done"
